$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 0.3333333333333333
$ws.Range("G2").Value = 0.002837
$ws.Range("H2").Value = 0.008510999999999999
$ws.Range("I2").Value = 0.00007108247730492929
$ws.Range("J2").Value = 0.00007108247730492929
$ws.Range("M2").Value = 16.27546433333333
$ws.Range("N2").Value = 48.826393
$ws.Range("O2").Value = 0.06628560529319844
$ws.Range("P2").Value = 0.06628560529319844
$ws.Range("Q2").Value = 0.04617349231366666
$ws.Range("R2").Value = 0.415561430823
$ws.Range("S2").Value = 0.000004711745033897278
$ws.Range("T2").Value = 0.000004711745033897278
# Row 3
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 0.3333333333333333
$ws.Range("G3").Value = 0.002837
$ws.Range("H3").Value = 0.008510999999999999
$ws.Range("I3").Value = 0.00007108247730492929
$ws.Range("J3").Value = 0.00007108247730492929
$ws.Range("O3").Value = 0.3480686258826592
$ws.Range("P3").Value = 0.3480686258826592
$ws.Range("Q3").Value = 0.2424590369316667
$ws.Range("R3").Value = 2.182131332385
$ws.Range("S3").Value = 0.00002474158019986204
$ws.Range("T3").Value = 0.00002474158019986204
# Row 4
$ws.Range("E4").Value = 1
$ws.Range("F4").Value = 0.3333333333333333
$ws.Range("G4").Value = 0.002837
$ws.Range("H4").Value = 0.008510999999999999
$ws.Range("I4").Value = 0.00007108247730492929
$ws.Range("J4").Value = 0.00007108247730492929
$ws.Range("M4").Value = 42.61351133333333
$ws.Range("N4").Value = 127.840534
$ws.Range("O4").Value = 0.17355341356458
$ws.Range("P4").Value = 0.17355341356458
$ws.Range("Q4").Value = 0.1208945316526666
$ws.Range("R4").Value = 1.088050784874
$ws.Range("S4").Value = 0.00001233660658089726
$ws.Range("T4").Value = 0.00001233660658089726
# Row 5
$ws.Range("E5").Value = 1
$ws.Range("F5").Value = 0.3333333333333333
$ws.Range("G5").Value = 0.002837
$ws.Range("H5").Value = 0.008510999999999999
$ws.Range("I5").Value = 0.00007108247730492929
$ws.Range("J5").Value = 0.00007108247730492929
$ws.Range("M5").Value = 101.183272
$ws.Range("N5").Value = 303.549816
$ws.Range("O5").Value = 0.4120923552595624
$ws.Range("P5").Value = 0.4120923552595624
$ws.Range("Q5").Value = 0.2870569426639999
$ws.Range("R5").Value = 2.583512483976
$ws.Range("S5").Value = 0.0000292925454902727
$ws.Range("T5").Value = 0.0000292925454902727
# Row 6
$ws.Range("I6").Value = 0.3776915775490952
$ws.Range("J6").Value = 0.3776915775490952
$ws.Range("M6").Value = 16.27546433333333
$ws.Range("N6").Value = 48.826393
$ws.Range("O6").Value = 0.06628560529319844
$ws.Range("P6").Value = 0.06628560529319844
$ws.Range("Q6").Value = 245.3394959504378
$ws.Range("R6").Value = 2208.05546355394
$ws.Range("S6").Value = 0.02503551483198477
$ws.Range("T6").Value = 0.02503551483198477
# Row 7
$ws.Range("I7").Value = 0.3776915775490952
$ws.Range("J7").Value = 0.3776915775490952
$ws.Range("O7").Value = 0.3480686258826592
$ws.Range("P7").Value = 0.3480686258826592
$ws.Range("S7").Value = 0.1314625884049674
$ws.Range("T7").Value = 0.1314625884049674
# Row 8
$ws.Range("I8").Value = 0.3776915775490952
$ws.Range("J8").Value = 0.3776915775490952
$ws.Range("M8").Value = 42.61351133333333
$ws.Range("N8").Value = 127.840534
$ws.Range("O8").Value = 0.17355341356458
$ws.Range("P8").Value = 0.17355341356458
$ws.Range("Q8").Value = 642.3643084508577
$ws.Range("R8").Value = 5781.278776057719
$ws.Range("S8").Value = 0.06554966255823676
$ws.Range("T8").Value = 0.06554966255823677
# Row 9
$ws.Range("I9").Value = 0.3776915775490952
$ws.Range("J9").Value = 0.3776915775490952
$ws.Range("M9").Value = 101.183272
$ws.Range("N9").Value = 303.549816
$ws.Range("O9").Value = 0.4120923552595624
$ws.Range("P9").Value = 0.4120923552595624
$ws.Range("Q9").Value = 1525.256204227253
$ws.Range("R9").Value = 13727.30583804528
$ws.Range("S9").Value = 0.1556438117539063
$ws.Range("T9").Value = 0.1556438117539063
# Row 10
$ws.Range("G10").Value = 1.581618666666667
$ws.Range("H10").Value = 4.744856
$ws.Range("I10").Value = 0.03962825977384063
$ws.Range("J10").Value = 0.03962825977384063
$ws.Range("M10").Value = 16.27546433333333
$ws.Range("N10").Value = 48.826393
$ws.Range("O10").Value = 0.06628560529319844
$ws.Range("P10").Value = 0.06628560529319844
$ws.Range("Q10").Value = 25.74157819826755
$ws.Range("R10").Value = 231.674203784408
$ws.Range("S10").Value = 0.002626783185825133
$ws.Range("T10").Value = 0.002626783185825133
# Row 11
$ws.Range("G11").Value = 1.581618666666667
$ws.Range("H11").Value = 4.744856
$ws.Range("I11").Value = 0.03962825977384063
$ws.Range("J11").Value = 0.03962825977384063
$ws.Range("O11").Value = 0.3480686258826592
$ws.Range("P11").Value = 0.3480686258826592
$ws.Range("Q11").Value = 135.1701581646622
$ws.Range("R11").Value = 1216.53142348196
$ws.Range("S11").Value = 0.01379335392560177
$ws.Range("T11").Value = 0.01379335392560177
# Row 12
$ws.Range("G12").Value = 1.581618666666667
$ws.Range("H12").Value = 4.744856
$ws.Range("I12").Value = 0.03962825977384063
$ws.Range("J12").Value = 0.03962825977384063
$ws.Range("M12").Value = 42.61351133333333
$ws.Range("N12").Value = 127.840534
$ws.Range("O12").Value = 0.17355341356458
$ws.Range("P12").Value = 0.17355341356458
$ws.Range("Q12").Value = 67.39832497701155
$ws.Range("R12").Value = 606.5849247931041
$ws.Range("S12").Value = 0.006877619757373971
$ws.Range("T12").Value = 0.006877619757373972
# Row 13
$ws.Range("G13").Value = 1.581618666666667
$ws.Range("H13").Value = 4.744856
$ws.Range("I13").Value = 0.03962825977384063
$ws.Range("J13").Value = 0.03962825977384063
$ws.Range("M13").Value = 101.183272
$ws.Range("N13").Value = 303.549816
$ws.Range("O13").Value = 0.4120923552595624
$ws.Range("P13").Value = 0.4120923552595624
$ws.Range("Q13").Value = 160.0333517496107
$ws.Range("R13").Value = 1440.300165746496
$ws.Range("S13").Value = 0.01633050290503976
$ws.Range("T13").Value = 0.01633050290503976
# Row 14
$ws.Range("G14").Value = 23.25273433333334
$ws.Range("H14").Value = 69.75820300000001
$ws.Range("I14").Value = 0.5826090801997593
$ws.Range("J14").Value = 0.5826090801997593
$ws.Range("M14").Value = 16.27546433333333
$ws.Range("N14").Value = 48.826393
$ws.Range("O14").Value = 0.06628560529319844
$ws.Range("P14").Value = 0.06628560529319844
$ws.Range("Q14").Value = 378.4490482946421
$ws.Range("R14").Value = 3406.041434651779
$ws.Range("S14").Value = 0.03861859553035464
$ws.Range("T14").Value = 0.03861859553035464
# Row 15
$ws.Range("G15").Value = 23.25273433333334
$ws.Range("H15").Value = 69.75820300000001
$ws.Range("I15").Value = 0.5826090801997593
$ws.Range("J15").Value = 0.5826090801997593
$ws.Range("O15").Value = 0.3480686258826592
$ws.Range("P15").Value = 0.3480686258826592
$ws.Range("Q15").Value = 1987.25258106729
$ws.Range("R15").Value = 17885.27322960561
$ws.Range("S15").Value = 0.2027879419718902
$ws.Range("T15").Value = 0.2027879419718902
# Row 16
$ws.Range("G16").Value = 23.25273433333334
$ws.Range("H16").Value = 69.75820300000001
$ws.Range("I16").Value = 0.5826090801997593
$ws.Range("J16").Value = 0.5826090801997593
$ws.Range("M16").Value = 42.61351133333333
$ws.Range("N16").Value = 127.840534
$ws.Range("O16").Value = 0.17355341356458
$ws.Range("P16").Value = 0.17355341356458
$ws.Range("Q16").Value = 990.8806580444891
$ws.Range("R16").Value = 8917.925922400402
$ws.Range("S16").Value = 0.1011137946423884
$ws.Range("T16").Value = 0.1011137946423884
# Row 17
$ws.Range("G17").Value = 23.25273433333334
$ws.Range("H17").Value = 69.75820300000001
$ws.Range("I17").Value = 0.5826090801997593
$ws.Range("J17").Value = 0.5826090801997593
$ws.Range("M17").Value = 101.183272
$ws.Range("N17").Value = 303.549816
$ws.Range("O17").Value = 0.4120923552595624
$ws.Range("P17").Value = 0.4120923552595624
$ws.Range("Q17").Value = 2352.787742793405
$ws.Range("R17").Value = 21175.08968514065
$ws.Range("S17").Value = 0.2400887480551261
$ws.Range("T17").Value = 0.2400887480551261
